$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so values like
# "325.45" or "1.020" are not auto-converted to numbers (which would
# lose trailing zeros / introduce floating point artifacts).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.497.36'
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.921.01'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.45'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4839'
$ws.Range("E7").Value = '  +3.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4086'
$ws.Range("E8").Value = '  +1.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08173'
$ws.Range("E9").Value = '  +2.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.020'
$ws.Range("E10").Value = '  +3.04%  '
$ws.Range("E11").Value = '  +5.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.926.67'
$ws.Range("E12").Value = '  +2.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.053'
$ws.Range("E13").Value = '  +3.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.238'
$ws.Range("E14").Value = '  +2.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.42'
$ws.Range("E15").Value = '  +3.11%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06750'
$ws.Range("E17").Value = '  +2.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001039'
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.79'
$ws.Range("E19").Value = '  +2.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.526.65'
$ws.Range("E21").Value = '  +1.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.637'
$ws.Range("E22").Value = '  +2.73%  '
$ws.Range("E23").Value = '  +2.31%  '
$ws.Range("E24").Value = '  -0.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.155.98'
$ws.Range("E25").Value = '  +1.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.688'
$ws.Range("E26").Value = '  +12.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.75'
$ws.Range("E27").Value = '  +1.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.08'
$ws.Range("E28").Value = '  +2.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.129'
$ws.Range("E29").Value = '  +2.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.68'
$ws.Range("E30").Value = '  +3.22%  '
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09564'
$ws.Range("E32").Value = '  +1.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.544'
$ws.Range("E33").Value = '  +3.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.398'
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.559'
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02284'
$ws.Range("E36").Value = '  +2.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06138'
$ws.Range("E37").Value = '  +1.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.187'
$ws.Range("E38").Value = '  +1.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.87'
$ws.Range("E39").Value = '  +8.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5984'
$ws.Range("E40").Value = '  +3.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.019'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1862'
$ws.Range("E42").Value = '  +2.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.429'
$ws.Range("E43").Value = '  -1.04%  '
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07631'
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.46'
$ws.Range("E46").Value = '  +2.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5588'
$ws.Range("E47").Value = '  +2.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.961'
$ws.Range("E48").Value = '  +3.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.80'
$ws.Range("E49").Value = '  +3.15%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.435'
$ws.Range("E50").Value = '  +4.39%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.88'
$ws.Range("E51").Value = '  +3.14%  '
